# Osteopetrosis.xlsx: refresh "data" timestamps and add a "metadata" tab
# summarising the PanelApp query that produced this export.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("data")

# ---------------------------------------------------------------------
# 1. Refresh the per-row "time_taken" timestamps on the "data" sheet.
# ---------------------------------------------------------------------
$data.Range("F2").Value  = "2021-10-05 14:21:56.179056"
$data.Range("F3").Value  = "2021-10-05 14:21:56.179063"
$data.Range("F4").Value  = "2021-10-05 14:21:56.179066"
$data.Range("F5").Value  = "2021-10-05 14:21:56.179069"
$data.Range("F6").Value  = "2021-10-05 14:21:56.179072"
$data.Range("F7").Value  = "2021-10-05 14:21:56.179074"
$data.Range("F8").Value  = "2021-10-05 14:21:56.179077"
$data.Range("F9").Value  = "2021-10-05 14:21:56.179079"
$data.Range("F10").Value = "2021-10-05 14:21:56.179082"
$data.Range("F11").Value = "2021-10-05 14:21:56.179085"
$data.Range("F12").Value = "2021-10-05 14:21:56.179087"
$data.Range("F13").Value = "2021-10-05 14:21:56.179089"
$data.Range("F14").Value = "2021-10-05 14:21:56.179092"
$data.Range("F15").Value = "2021-10-05 14:21:56.179095"
$data.Range("F16").Value = "2021-10-05 14:21:56.179097"
$data.Range("F17").Value = "2021-10-05 14:21:56.179100"
$data.Range("F18").Value = "2021-10-05 14:21:56.179103"
$data.Range("F19").Value = "2021-10-05 14:21:56.179105"
$data.Range("F20").Value = "2021-10-05 14:21:56.179108"
$data.Range("F21").Value = "2021-10-05 14:21:56.179110"

# ---------------------------------------------------------------------
# 2. Add a new "metadata" worksheet right after "data".
# ---------------------------------------------------------------------
$meta = $wb.Worksheets.Add($null, $data)
$meta.Name = "metadata"

# Header row (bold / centered / bordered, matching the "data" sheet's
# header style) — copy formatting from "data" so the same style index
# is reused instead of a new (merely equivalent) one being created.
$data.Range("B1:F1").Copy()
$meta.Range("B1:F1").PasteSpecial(-4122)
$data.Range("F1").Copy()
$meta.Range("G1").PasteSpecial(-4122)

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

# Row 2 index cell (A2) uses the same bold/bordered style as "data"!A2.
$data.Range("A2").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$meta.Range("A2").Value = 0

$meta.Range("B2").Value = "Osteopetrosis"
$meta.Range("C2").Value = 943

# "1.26" must stay a text value (not be coerced to the number 1.26) and
# must keep the sheet's default (unstyled) cell format. Stage it via a
# scratch cell so the text-number-format xf never attaches to D2 itself.
$meta.Range("Z1").NumberFormat = "@"
$meta.Range("Z1").Value = "1.26"
$meta.Range("Z1").Copy()
$meta.Range("D2").PasteSpecial(-4163)
$meta.Range("Z1").Clear()

$meta.Range("E2").Value = "2021-03-23T22:34:38.909988Z"
$meta.Range("F2").Value = "2021-10-05 14:21:56.175474"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/943/?format=json"

$meta.Range("A1").Select() | Out-Null

# Keep "data" as the active sheet (its bookView wasn't touched by the
# commit — only a new <sheet> entry was appended).
$data.Activate() | Out-Null
